$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 6, pushing existing data (rows 6..68) down to 7..69.
$ws.Rows("6:6").Insert()

# Populate the newly inserted row 6 with the new data record.
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Macroferia Regional de Talca"
$ws.Range("C6").Value = "Maule"
$ws.Range("D6").Value = 44490
$ws.Range("E6").Value = 7
$ws.Range("F6").Value = 100112013
$ws.Range("G6").Value = "Alcachofa"
$ws.Range("H6").Value = "Madrigal"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 300
$ws.Range("K6").Value = 9000
$ws.Range("L6").Value = 9000
$ws.Range("M6").Value = 9000
$ws.Range("N6").Value = "`$/caja 40 unidades"
$ws.Range("O6").Value = "Provincia del Elquí"
$ws.Range("P6").Value = 225
$ws.Range("Q6").Value = 40
$ws.Range("R6").Value = "Hortaliza"
